$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (column D). Cells D4,D5,D6,D7,D8,D12,D13 all share the
# same underlying string value, so update every one of them to keep them sharing
# a single (updated) shared-string entry, same as the source diff.
$ws.Range("D4").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D5").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D6").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D7").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D8").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D12").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("D13").Value = "2024-08-12T02:00:00.000Z"

# Update numeric properties in row 7
$ws.Range("W7").Value = 48813000
$ws.Range("AA7").Value = 45937000
$ws.Range("AE7").Value = 94750000
$ws.Range("AH7").Value = 75250000
$ws.Range("AK7").Value = 13
$ws.Range("AN7").Value = 19500000
$ws.Range("AQ7").Value = 85250000
